# Auto-generated edit script applying the row-data rotation / odds update
# described by the commit diff (Greece Super League 1 sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74: full B:AC content rotates in from another fixture row
$ws.Range("B74").Value2 = 5369548
$ws.Range("F74").Value2 = "Atromitos Athinon"
$ws.Range("G74").Value2 = "Panathinaikos"
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 2
$ws.Range("J74").Value2 = "A"
$ws.Range("K74").Value2 = 5.75
$ws.Range("L74").Value2 = 3.25
$ws.Range("M74").Value2 = 1.75
$ws.Range("N74").Value2 = 7
$ws.Range("O74").Value2 = 3.6
$ws.Range("P74").Value2 = 1.571
$ws.Range("Q74").Value2 = 0.75
$ws.Range("R74").Value2 = 2.05
$ws.Range("S74").Value2 = 1.75
$ws.Range("T74").Value2 = 2.25
$ws.Range("U74").Value2 = 1.95
$ws.Range("V74").Value2 = 1.9
$ws.Range("W74").Value2 = -1
$ws.Range("X74").Value2 = -1
$ws.Range("Y74").Value2 = 0.571
$ws.Range("Z74").Value2 = -1
$ws.Range("AA74").Value2 = 0.75
$ws.Range("AB74").Value2 = -0.5
$ws.Range("AC74").Value2 = 0.45

# Row 75: full B:AC content rotates in from another fixture row
$ws.Range("B75").Value2 = 5369519
$ws.Range("F75").Value2 = "Levadiakos"
$ws.Range("G75").Value2 = "OFI Crete"
$ws.Range("H75").Value2 = 2
$ws.Range("I75").Value2 = 0
$ws.Range("J75").Value2 = "H"
$ws.Range("K75").Value2 = 2.55
$ws.Range("L75").Value2 = 3.25
$ws.Range("M75").Value2 = 2.8
$ws.Range("N75").Value2 = 2.9
$ws.Range("O75").Value2 = 3.3
$ws.Range("P75").Value2 = 2.45
$ws.Range("Q75").Value2 = 0
$ws.Range("R75").Value2 = 2.075
$ws.Range("S75").Value2 = 1.725
$ws.Range("T75").Value2 = 2
$ws.Range("U75").Value2 = 1.95
$ws.Range("V75").Value2 = 1.9
$ws.Range("W75").Value2 = 1.9
$ws.Range("X75").Value2 = -1
$ws.Range("Y75").Value2 = -1
$ws.Range("Z75").Value2 = 1.075
$ws.Range("AA75").Value2 = -1
$ws.Range("AB75").Value2 = 0
$ws.Range("AC75").Value2 = -0

# Row 76: full B:AC content rotates in from another fixture row
$ws.Range("B76").Value2 = 5374242
$ws.Range("F76").Value2 = "Volos NFC"
$ws.Range("G76").Value2 = "PAOK Salonika"
$ws.Range("H76").Value2 = 0
$ws.Range("I76").Value2 = 1
$ws.Range("J76").Value2 = "A"
$ws.Range("K76").Value2 = 12
$ws.Range("L76").Value2 = 5
$ws.Range("M76").Value2 = 1.3
$ws.Range("N76").Value2 = 11
$ws.Range("O76").Value2 = 4.75
$ws.Range("P76").Value2 = 1.285
$ws.Range("Q76").Value2 = 1.5
$ws.Range("R76").Value2 = 1.875
$ws.Range("S76").Value2 = 1.975
$ws.Range("T76").Value2 = 2.5
$ws.Range("U76").Value2 = 1.875
$ws.Range("V76").Value2 = 1.975
$ws.Range("W76").Value2 = -1
$ws.Range("X76").Value2 = -1
$ws.Range("Y76").Value2 = 0.2849999999999999
$ws.Range("Z76").Value2 = 0.875
$ws.Range("AA76").Value2 = -1
$ws.Range("AB76").Value2 = -1
$ws.Range("AC76").Value2 = 0.9750000000000001

# Row 77: full B:AC content rotates in from another fixture row
$ws.Range("B77").Value2 = 5374241
$ws.Range("F77").Value2 = "Aris Salonika"
$ws.Range("G77").Value2 = "Giannina"
$ws.Range("H77").Value2 = 3
$ws.Range("I77").Value2 = 1
$ws.Range("J77").Value2 = "H"
$ws.Range("K77").Value2 = 1.4
$ws.Range("L77").Value2 = 4.333
$ws.Range("M77").Value2 = 9
$ws.Range("N77").Value2 = 1.444
$ws.Range("O77").Value2 = 4
$ws.Range("P77").Value2 = 9
$ws.Range("Q77").Value2 = -1.25
$ws.Range("R77").Value2 = 1.975
$ws.Range("S77").Value2 = 1.875
$ws.Range("T77").Value2 = 2.25
$ws.Range("U77").Value2 = 1.85
$ws.Range("V77").Value2 = 2
$ws.Range("W77").Value2 = 0.444
$ws.Range("X77").Value2 = -1
$ws.Range("Y77").Value2 = -1
$ws.Range("Z77").Value2 = 0.9750000000000001
$ws.Range("AA77").Value2 = -1
$ws.Range("AB77").Value2 = 0.8500000000000001
$ws.Range("AC77").Value2 = -1

# Row 78: full B:AC content rotates in from another fixture row
$ws.Range("B78").Value2 = 5374240
$ws.Range("F78").Value2 = "AEK Athens"
$ws.Range("G78").Value2 = "Olympiakos"
$ws.Range("H78").Value2 = 1
$ws.Range("I78").Value2 = 3
$ws.Range("J78").Value2 = "A"
$ws.Range("K78").Value2 = 2.2
$ws.Range("L78").Value2 = 3.25
$ws.Range("M78").Value2 = 3.4
$ws.Range("N78").Value2 = 1.85
$ws.Range("O78").Value2 = 3.4
$ws.Range("P78").Value2 = 4.333
$ws.Range("Q78").Value2 = -0.5
$ws.Range("R78").Value2 = 1.875
$ws.Range("S78").Value2 = 1.975
$ws.Range("T78").Value2 = 2.5
$ws.Range("U78").Value2 = 2
$ws.Range("V78").Value2 = 1.85
$ws.Range("W78").Value2 = -1
$ws.Range("X78").Value2 = -1
$ws.Range("Y78").Value2 = 3.333
$ws.Range("Z78").Value2 = -1
$ws.Range("AA78").Value2 = 0.9750000000000001
$ws.Range("AB78").Value2 = 1
$ws.Range("AC78").Value2 = -1

# Row 103: full B:AC content rotates in from another fixture row
$ws.Range("B103").Value2 = 6399628
$ws.Range("F103").Value2 = "Lamia"
$ws.Range("G103").Value2 = "Atromitos Athinon"
$ws.Range("H103").Value2 = 1
$ws.Range("I103").Value2 = 0
$ws.Range("J103").Value2 = "H"
$ws.Range("K103").Value2 = 2.2
$ws.Range("L103").Value2 = 3.2
$ws.Range("M103").Value2 = 3.5
$ws.Range("N103").Value2 = 1.8
$ws.Range("O103").Value2 = 3.5
$ws.Range("P103").Value2 = 4.75
$ws.Range("Q103").Value2 = -0.5
$ws.Range("R103").Value2 = 1.825
$ws.Range("S103").Value2 = 2.025
$ws.Range("T103").Value2 = 2.25
$ws.Range("U103").Value2 = 2.05
$ws.Range("V103").Value2 = 1.8
$ws.Range("W103").Value2 = 0.8
$ws.Range("X103").Value2 = -1
$ws.Range("Y103").Value2 = -1
$ws.Range("Z103").Value2 = 0.825
$ws.Range("AA103").Value2 = -1
$ws.Range("AB103").Value2 = -1
$ws.Range("AC103").Value2 = 0.8

# Row 104: full B:AC content rotates in from another fixture row
$ws.Range("B104").Value2 = 6399629
$ws.Range("F104").Value2 = "Levadiakos"
$ws.Range("G104").Value2 = "Ionikos Nikea"
$ws.Range("H104").Value2 = 2
$ws.Range("I104").Value2 = 2
$ws.Range("J104").Value2 = "D"
$ws.Range("K104").Value2 = 2.4
$ws.Range("L104").Value2 = 3.25
$ws.Range("M104").Value2 = 3
$ws.Range("N104").Value2 = 2.25
$ws.Range("O104").Value2 = 3.1
$ws.Range("P104").Value2 = 3.5
$ws.Range("Q104").Value2 = -0.25
$ws.Range("R104").Value2 = 1.9
$ws.Range("S104").Value2 = 1.95
$ws.Range("T104").Value2 = 1.75
$ws.Range("U104").Value2 = 1.825
$ws.Range("V104").Value2 = 2.025
$ws.Range("W104").Value2 = -1
$ws.Range("X104").Value2 = 2.1
$ws.Range("Y104").Value2 = -1
$ws.Range("Z104").Value2 = -0.5
$ws.Range("AA104").Value2 = 0.475
$ws.Range("AB104").Value2 = 0.825
$ws.Range("AC104").Value2 = -1

# Row 105: full B:AC content rotates in from another fixture row
$ws.Range("B105").Value2 = 6399627
$ws.Range("F105").Value2 = "Giannina"
$ws.Range("G105").Value2 = "Panetolikos"
$ws.Range("H105").Value2 = 3
$ws.Range("I105").Value2 = 2
$ws.Range("J105").Value2 = "H"
$ws.Range("K105").Value2 = 2.2
$ws.Range("L105").Value2 = 3.1
$ws.Range("M105").Value2 = 3.6
$ws.Range("N105").Value2 = 1.8
$ws.Range("O105").Value2 = 3.1
$ws.Range("P105").Value2 = 5.75
$ws.Range("Q105").Value2 = -0.75
$ws.Range("R105").Value2 = 2.05
$ws.Range("S105").Value2 = 1.8
$ws.Range("T105").Value2 = 2
$ws.Range("U105").Value2 = 1.825
$ws.Range("V105").Value2 = 2.025
$ws.Range("W105").Value2 = 0.8
$ws.Range("X105").Value2 = -1
$ws.Range("Y105").Value2 = -1
$ws.Range("Z105").Value2 = 0.5249999999999999
$ws.Range("AA105").Value2 = -0.5
$ws.Range("AB105").Value2 = 0.825
$ws.Range("AC105").Value2 = -1

# Row 128: full B:AC content rotates in from another fixture row
$ws.Range("B128").Value2 = 6397969
$ws.Range("F128").Value2 = "Olympiakos"
$ws.Range("G128").Value2 = "Panathinaikos"
$ws.Range("H128").Value2 = 1
$ws.Range("I128").Value2 = 0
$ws.Range("J128").Value2 = "H"
$ws.Range("K128").Value2 = 2.45
$ws.Range("L128").Value2 = 3.25
$ws.Range("M128").Value2 = 2.7
$ws.Range("N128").Value2 = 3.2
$ws.Range("O128").Value2 = 3.25
$ws.Range("P128").Value2 = 2.375
$ws.Range("Q128").Value2 = 0.25
$ws.Range("R128").Value2 = 1.825
$ws.Range("S128").Value2 = 2.025
$ws.Range("T128").Value2 = 2
$ws.Range("U128").Value2 = 1.75
$ws.Range("V128").Value2 = 2.05
$ws.Range("W128").Value2 = 2.2
$ws.Range("X128").Value2 = -1
$ws.Range("Y128").Value2 = -1
$ws.Range("Z128").Value2 = 0.825
$ws.Range("AA128").Value2 = -1
$ws.Range("AB128").Value2 = -1
$ws.Range("AC128").Value2 = 1.05

# Row 129: full B:AC content rotates in from another fixture row
$ws.Range("B129").Value2 = 6399615
$ws.Range("F129").Value2 = "Aris Salonika"
$ws.Range("G129").Value2 = "AEK Athens"
$ws.Range("H129").Value2 = 1
$ws.Range("I129").Value2 = 2
$ws.Range("J129").Value2 = "A"
$ws.Range("K129").Value2 = 4.5
$ws.Range("L129").Value2 = 3.25
$ws.Range("M129").Value2 = 1.8
$ws.Range("N129").Value2 = 12
$ws.Range("O129").Value2 = 5.25
$ws.Range("P129").Value2 = 1.285
$ws.Range("Q129").Value2 = 1.5
$ws.Range("R129").Value2 = 2
$ws.Range("S129").Value2 = 1.85
$ws.Range("T129").Value2 = 2.75
$ws.Range("U129").Value2 = 1.875
$ws.Range("V129").Value2 = 1.975
$ws.Range("W129").Value2 = -1
$ws.Range("X129").Value2 = -1
$ws.Range("Y129").Value2 = 0.2849999999999999
$ws.Range("Z129").Value2 = 1
$ws.Range("AA129").Value2 = -1
$ws.Range("AB129").Value2 = 0.4375
$ws.Range("AC129").Value2 = -0.5

# Row 259: full B:AC content rotates in from another fixture row
$ws.Range("B259").Value2 = 6936857
$ws.Range("F259").Value2 = "AEK Athens"
$ws.Range("G259").Value2 = "Panathinaikos"
$ws.Range("H259").Value2 = 2
$ws.Range("I259").Value2 = 2
$ws.Range("J259").Value2 = "D"
$ws.Range("K259").Value2 = 1.909
$ws.Range("L259").Value2 = 3.5
$ws.Range("M259").Value2 = 4.2
$ws.Range("N259").Value2 = 2.15
$ws.Range("O259").Value2 = 3.2
$ws.Range("P259").Value2 = 3.5
$ws.Range("Q259").Value2 = -0.25
$ws.Range("R259").Value2 = 1.85
$ws.Range("S259").Value2 = 2
$ws.Range("T259").Value2 = 2
$ws.Range("U259").Value2 = 1.8
$ws.Range("V259").Value2 = 2.05
$ws.Range("W259").Value2 = -1
$ws.Range("X259").Value2 = 2.2
$ws.Range("Y259").Value2 = -1
$ws.Range("Z259").Value2 = -0.5
$ws.Range("AA259").Value2 = 0.5
$ws.Range("AB259").Value2 = 0.8
$ws.Range("AC259").Value2 = -1

# Row 260: full B:AC content rotates in from another fixture row
$ws.Range("B260").Value2 = 6937238
$ws.Range("F260").Value2 = "PAOK Salonika"
$ws.Range("G260").Value2 = "Giannina"
$ws.Range("H260").Value2 = 4
$ws.Range("I260").Value2 = 0
$ws.Range("J260").Value2 = "H"
$ws.Range("K260").Value2 = 1.111
$ws.Range("L260").Value2 = 9
$ws.Range("M260").Value2 = 23
$ws.Range("N260").Value2 = 1.25
$ws.Range("O260").Value2 = 6
$ws.Range("P260").Value2 = 9
$ws.Range("Q260").Value2 = -1.75
$ws.Range("R260").Value2 = 2.025
$ws.Range("S260").Value2 = 1.825
$ws.Range("T260").Value2 = 2.75
$ws.Range("U260").Value2 = 1.8
$ws.Range("V260").Value2 = 2.05
$ws.Range("W260").Value2 = 0.25
$ws.Range("X260").Value2 = -1
$ws.Range("Y260").Value2 = -1
$ws.Range("Z260").Value2 = 1.025
$ws.Range("AA260").Value2 = -1
$ws.Range("AB260").Value2 = 0.8
$ws.Range("AC260").Value2 = -1

# Row 296: odds update (future fixture, only pre-match odds columns change)
$ws.Range("N296").Value2 = 9
$ws.Range("O296").Value2 = 4.75
$ws.Range("P296").Value2 = 1.363
$ws.Range("Q296").Value2 = 1.5
$ws.Range("R296").Value2 = 1.85
$ws.Range("S296").Value2 = 2
$ws.Range("U296").Value2 = 1.85
$ws.Range("V296").Value2 = 2

# Row 297: odds update (future fixture, only pre-match odds columns change)
$ws.Range("N297").Value2 = 2
$ws.Range("O297").Value2 = 3.4
$ws.Range("P297").Value2 = 3.8
$ws.Range("Q297").Value2 = -0.5
$ws.Range("U297").Value2 = 1.825
$ws.Range("V297").Value2 = 2.025

# Row 298: odds update (future fixture, only pre-match odds columns change)
$ws.Range("U298").Value2 = 1.95
$ws.Range("V298").Value2 = 1.9

# Row 299: odds update (future fixture, only pre-match odds columns change)
$ws.Range("N299").Value2 = 7.5
$ws.Range("Q299").Value2 = 1.25
$ws.Range("R299").Value2 = 2
$ws.Range("S299").Value2 = 1.85
$ws.Range("U299").Value2 = 1.875
$ws.Range("V299").Value2 = 1.975

# Row 300: odds update (future fixture, only pre-match odds columns change)
$ws.Range("R300").Value2 = 1.8
$ws.Range("S300").Value2 = 2.05

# Row 301: odds update (future fixture, only pre-match odds columns change)
$ws.Range("N301").Value2 = 1.5
$ws.Range("P301").Value2 = 7
$ws.Range("R301").Value2 = 1.825
$ws.Range("S301").Value2 = 2.025

# Row 302: odds update (future fixture, only pre-match odds columns change)
$ws.Range("R302").Value2 = 1.925
$ws.Range("S302").Value2 = 1.925
